# Auto-generated Excel COM-interop script to apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.906.68"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.402.14"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'561.85"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'142.21"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.530"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'5.26"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.76%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'25.50"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.25%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.0000172"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.838.42"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'61.828.78"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.408.73"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'11.21"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'321.18"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'6.80"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'66.12"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'1.75"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'8.78"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.71%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'562.50"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.66%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.520.28"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'0.0₃0931"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'8.16"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.76%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.39"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.77%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.146"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.51"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'4.66"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.50%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'5.42"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.09%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'Monero"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'151.27"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'PolygonEcosystemToken"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'0.378"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'18.56"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '  -8.19%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'147.35"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'3.60"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.0529"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'19.80"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.586"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.0915"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0224"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '  +0.59%  '
$ws.Range('E51').Style = 'Normal'
